# Auto-generated edit script: updates the cryptos list cell values per the
# "Updated cryptos list on Mon Nov 25 10:39:53 UTC 2024 with GitHub Actions"
# commit. All touched cells keep their original plain-text (inline/shared
# string) storage -- column D holds price strings that are NOT valid
# locale-neutral numbers (e.g. "98.449.59"), so every cell is written the
# same way regardless of whether it happens to parse as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values must stay text even when they look numeric
# (e.g. "253.96", "9.00", "0.0000260") -- otherwise Excel's smart entry
# would silently coerce them to doubles and mangle the trailing zeros /
# drop the text formatting entirely. Forcing the number format to Text
# before the write, then resetting the cell Style back to Normal afterwards,
# keeps the value as a literal string without leaving any residual styling
# behind on the cell.
$priceCells = @(
    @{Cell='D2'; Value='98.449.59'}
    @{Cell='D3'; Value='3.491.86'}
    @{Cell='D5'; Value='253.96'}
    @{Cell='D6'; Value='667.48'}
    @{Cell='D7'; Value='1.49'}
    @{Cell='D8'; Value='0.429'}
    @{Cell='D11'; Value='3.485.36'}
    @{Cell='D12'; Value='45.99'}
    @{Cell='D13'; Value='0.210'}
    @{Cell='D14'; Value='98.359.97'}
    @{Cell='D15'; Value='6.18'}
    @{Cell='D16'; Value='4.158.91'}
    @{Cell='D17'; Value='0.0000260'}
    @{Cell='D18'; Value='9.00'}
    @{Cell='D19'; Value='3.499.33'}
    @{Cell='D20'; Value='18.76'}
    @{Cell='D21'; Value='11.76'}
    @{Cell='D22'; Value='0.530'}
    @{Cell='D23'; Value='521.02'}
    @{Cell='D24'; Value='3.44'}
    @{Cell='D26'; Value='6.79'}
    @{Cell='D27'; Value='97.74'}
    @{Cell='D28'; Value='12.63'}
    @{Cell='D29'; Value='12.42'}
    @{Cell='D30'; Value='2.91'}
    @{Cell='D31'; Value='0.146'}
    @{Cell='D34'; Value='0.587'}
    @{Cell='D36'; Value='30.60'}
    @{Cell='D37'; Value='1.53'}
    @{Cell='D38'; Value='8.02'}
    @{Cell='D39'; Value='0.155'}
    @{Cell='D40'; Value='530.11'}
    @{Cell='D42'; Value='0.907'}
    @{Cell='D43'; Value='1.78'}
    @{Cell='D44'; Value='24.43'}
    @{Cell='D45'; Value='0.0432'}
    @{Cell='D46'; Value='5.77'}
    @{Cell='D47'; Value='3.66'}
    @{Cell='D48'; Value='8.67'}
    @{Cell='D49'; Value='2.22'}
    @{Cell='D50'; Value='55.26'}
    @{Cell='D51'; Value='3.25'}
)

foreach ($u in $priceCells) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}

# Columns B (Coin), C (Link) and E (Volume(1h)) can be written directly --
# none of their new values are numeric-looking, so Excel keeps them as text.
$textCells = @(
    @{Cell='E2'; Value='  +0.68%  '}
    @{Cell='E3'; Value='  +2.92%  '}
    @{Cell='E4'; Value='  -0.12%  '}
    @{Cell='E5'; Value='  +0.24%  '}
    @{Cell='E6'; Value='  -0.65%  '}
    @{Cell='E7'; Value='  +3.30%  '}
    @{Cell='E8'; Value='  +0.27%  '}
    @{Cell='E9'; Value='  +1.41%  '}
    @{Cell='E10'; Value='  -0.06%  '}
    @{Cell='E11'; Value='  +2.82%  '}
    @{Cell='E12'; Value='  +10.79%  '}
    @{Cell='E13'; Value='  -1.42%  '}
    @{Cell='E14'; Value='  +0.74%  '}
    @{Cell='E15'; Value='  -0.64%  '}
    @{Cell='B16'; Value='WrappedliquidstakedEther2.0'}
    @{Cell='C16'; Value='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'}
    @{Cell='E16'; Value='  +3.25%  '}
    @{Cell='B17'; Value='ShibaInu'}
    @{Cell='C17'; Value='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'}
    @{Cell='E17'; Value='  -0.78%  '}
    @{Cell='E18'; Value='  +2.46%  '}
    @{Cell='E19'; Value='  +2.72%  '}
    @{Cell='E20'; Value='  +8.63%  '}
    @{Cell='B21'; Value='Uniswap'}
    @{Cell='C21'; Value='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'}
    @{Cell='E21'; Value='  +6.46%  '}
    @{Cell='B22'; Value='Stellar'}
    @{Cell='C22'; Value='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'}
    @{Cell='E22'; Value='  -8.39%  '}
    @{Cell='E23'; Value='  +2.53%  '}
    @{Cell='E24'; Value='  +0.72%  '}
    @{Cell='E25'; Value='  +0.47%  '}
    @{Cell='E26'; Value='  +5.08%  '}
    @{Cell='E27'; Value='  -1.98%  '}
    @{Cell='E28'; Value='  +1.53%  '}
    @{Cell='E29'; Value='  +8.93%  '}
    @{Cell='B30'; Value='PancakeSwap'}
    @{Cell='C30'; Value='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'}
    @{Cell='E30'; Value='  +10.53%  '}
    @{Cell='B31'; Value='Hedera'}
    @{Cell='C31'; Value='https://coinranking.com/coin/jad286TjB+hedera-hbar'}
    @{Cell='E31'; Value='  -2.96%  '}
    @{Cell='E33'; Value='  -1.26%  '}
    @{Cell='E34'; Value='  +4.02%  '}
    @{Cell='E35'; Value='  +0.31%  '}
    @{Cell='E36'; Value='  +4.78%  '}
    @{Cell='B37'; Value='Fetch.AI'}
    @{Cell='C37'; Value='https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'}
    @{Cell='E37'; Value='  +1.78%  '}
    @{Cell='B38'; Value='RenderToken'}
    @{Cell='C38'; Value='https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'}
    @{Cell='E38'; Value='  +0.86%  '}
    @{Cell='E39'; Value='  +2.79%  '}
    @{Cell='E40'; Value='  -1.33%  '}
    @{Cell='E41'; Value='  -0.11%  '}
    @{Cell='E42'; Value='  +5.11%  '}
    @{Cell='E43'; Value='  +4.77%  '}
    @{Cell='B44'; Value='WhiteBITCoin'}
    @{Cell='C44'; Value='https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'}
    @{Cell='E44'; Value='  -1.10%  '}
    @{Cell='B45'; Value='VeChain'}
    @{Cell='C45'; Value='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'}
    @{Cell='E45'; Value='  -0.20%  '}
    @{Cell='E46'; Value='  +0.09%  '}
    @{Cell='B47'; Value='MantraDAO'}
    @{Cell='C47'; Value='https://coinranking.com/coin/cTdD8lD-6+mantradao-om'}
    @{Cell='E47'; Value='  -1.53%  '}
    @{Cell='B48'; Value='Cosmos'}
    @{Cell='C48'; Value='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'}
    @{Cell='E48'; Value='  -4.66%  '}
    @{Cell='E49'; Value='  +7.97%  '}
    @{Cell='E50'; Value='  +1.63%  '}
    @{Cell='E51'; Value='  +1.66%  '}
)

foreach ($u in $textCells) {
    $ws.Range($u.Cell).Value = $u.Value
}

Write-Output ("Applied " + ($priceCells.Count + $textCells.Count) + " cell updates.")
